# "Generate Report for Archive" — refresh localization-status report:
#  - flip the "Ready for handoff" status to "In Translation" everywhere it
#    appears (Overview summary columns + each language sheet's Status column)
#  - re-fit the now-narrower status columns to their content

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Overview sheet: zh-cn / de-de status columns (E, F) for both data rows
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# Per-language sheets: Status column (C) for both data rows
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# The shorter status text means those columns no longer need to be as wide;
# re-fit them to their (now shorter) contents.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
